# Modifieradmin.xlsx edit: "verification des séquentiels de admin"
# The nominal scenario step "4) Quideance vérifie que le formulaire est
# correctement renseigné" (row 8) is removed, and every subsequent
# numbered step / cross-reference in the use-case table is renumbered
# down by one (5->4, 6->5, 7->6; "point 5)"->"point 4)", "point 6)"->"point 5)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the whole row 8 - this shifts rows 9-12 up to 8-11 and
#    automatically drops the now-unreferenced shared string.
$ws.Rows(8).Delete()

# 2) Plain-text cells that just need their leading step number fixed.
$ws.Range("A8").Value = "4) Quideance vérifie que l'email soit rempli`nQuideance vérifie que le mot de passe soit rempli`nQuideance vérifie que le rôle est sélectionné"

$ws.Range("A9").Value = "5) Quideance enregistre le formulaire"

$ws.Range("A10").Value = "6) Quideance redirige vers la page d'aministration`nAffichage d'un message `"Modification effectuée avec succès`""

# 3) Rich-text cell B8: "A1) Le formulaire n'est pas renseigné correctement" (bold title)
#    followed by the normal-weight explanation paragraph, now referencing point 4).
$b8 = $ws.Range("B8")
$b8.Value = "A1) Le formulaire n'est pas renseigné correctement`na lieu au point 4) du nominal`nQuideance affiche un message d'erreur adéquat et retourne au point 3) du nominal`n"
$b8.Characters(1, 50).Font.Bold = $true
$b8.Characters(1, 50).Font.Size = 12
$b8.Characters(1, 50).Font.Color = 0
$b8.Characters(51, 112).Font.Bold = $false
$b8.Characters(51, 112).Font.Size = 11
$b8.Characters(51, 112).Font.Color = 0

# 4) Rich-text cell C9: "E2) L'enregistrement ne se valide pas" (bold title) +
#    normal explanation (now referencing point 5)) + bold red "exit" sentence.
$c9 = $ws.Range("C9")
$c9.Value = "E2) L'enregistrement ne se valide pas`na lieu au point 5) du nominal`nQuideance affiche un message d'erreur adéquat`nOn sort du USE CASE sur un échec"
$c9.Characters(1, 37).Font.Bold = $true
$c9.Characters(1, 37).Font.Size = 12
$c9.Characters(1, 37).Font.Color = 0
$c9.Characters(38, 77).Font.Bold = $false
$c9.Characters(38, 77).Font.Size = 11
$c9.Characters(38, 77).Font.Color = 0
$c9.Characters(115, 32).Font.Bold = $true
$c9.Characters(115, 32).Font.Size = 11
$c9.Characters(115, 32).Font.Color = 255

# 5) Update the view state to match where the edit left the selection.
$null = $ws.Range("C9").Select()

Write-Output "edit applied"
